# Generate Report for Handoff
# Adds two new handed-off files (933ab61c-... and bac28f38-...) as new rows
# to the Overview / zh-cn / de-de sheets, mirroring the shape of the
# existing "Ready for handoff" rows (row 5, af120165 pattern's 3dfced4f row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet -> new rows 6 and 7
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md"
$wsOverview.Range("B6").Value = "e2e\933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-12 16:51:29"

$wsOverview.Range("A7").Value = "bac28f38-e2d8-499d-92b0-c819b0177acc.md"
$wsOverview.Range("B7").Value = "e2e\bac28f38-e2d8-499d-92b0-c819b0177acc.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-12 16:51:29"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/933ab61c/e2e/933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md", "", "", "e2e\933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/bac28f38/e2e/bac28f38-e2d8-499d-92b0-c819b0177acc.md", "", "", "e2e\bac28f38-e2d8-499d-92b0-c819b0177acc.md")

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------
# zh-cn sheet -> new rows 6 and 7
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.236a43bf444b66486868f88a302d828b50f68feb.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-12 16:51:17"
$wsZhCn.Range("I6").Value = ""
$wsZhCn.Range("J6").Value = ""
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L6").Value = ""
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("N6").Value = ""
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Range("P6").Value = ""

$wsZhCn.Range("A7").Value = "bac28f38-e2d8-499d-92b0-c819b0177acc.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "bac28f38-e2d8-499d-92b0-c819b0177acc.f679037a0c5f04b2a7f2bbfa5d414662afbcf707.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-12 16:51:17"
$wsZhCn.Range("I7").Value = ""
$wsZhCn.Range("J7").Value = ""
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L7").Value = ""
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("N7").Value = ""
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Range("P7").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/933ab61c/e2e/933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md", "", "", "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/bac28f38/e2e/bac28f38-e2d8-499d-92b0-c819b0177acc.md", "", "", "bac28f38-e2d8-499d-92b0-c819b0177acc.md")

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------
# de-de sheet -> new rows 6 and 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.236a43bf444b66486868f88a302d828b50f68feb.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-12 16:51:29"
$wsDeDe.Range("I6").Value = ""
$wsDeDe.Range("J6").Value = ""
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L6").Value = ""
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("N6").Value = ""
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Range("P6").Value = ""

$wsDeDe.Range("A7").Value = "bac28f38-e2d8-499d-92b0-c819b0177acc.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "bac28f38-e2d8-499d-92b0-c819b0177acc.f679037a0c5f04b2a7f2bbfa5d414662afbcf707.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-12 16:51:29"
$wsDeDe.Range("I7").Value = ""
$wsDeDe.Range("J7").Value = ""
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L7").Value = ""
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("N7").Value = ""
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Range("P7").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/933ab61c/e2e/933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md", "", "", "933ab61c-ab9f-4e8b-9265-d6ce0782a7d7.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/bac28f38/e2e/bac28f38-e2d8-499d-92b0-c819b0177acc.md", "", "", "bac28f38-e2d8-499d-92b0-c819b0177acc.md")

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P7"))
